$wb = $excel.ActiveWorkbook

# The workbook used to carry 3 sheets that fed a multi-step "carga /
# preaplicacion / aplicacion final" pagaduria flow. This edit collapses
# everything into a single DataProvider-style sheet: "AplicacionPago".
[void]$wb.Worksheets.Item("PreaplicacionPagaduria").Delete()
[void]$wb.Worksheets.Item("AplicacionFinalPagaduria").Delete()

$sheet = $wb.Worksheets.Item("CarguePlanillaAlSistema")
$sheet.Name = "AplicacionPago"

# Start from a clean slate: wipe data AND formatting (incl. the old
# column-wide Text number format) so the new layout below isn't polluted
# by whatever the previous table had.
[void]$sheet.Cells.Clear()
foreach ($colIdx in 1..10) {
    [void]$sheet.Columns.Item($colIdx).ClearFormats()
}

# Header row
$sheet.Range("A1").Value = "IdPagaduria"
$sheet.Range("B1").Value = "Periodo"
$sheet.Range("C1").Value = "NombrePagaduria"
$sheet.Range("D1").Value = "RutaPagaduria"
$sheet.Range("E1").Value = "Ano"
$sheet.Range("F1").Value = "PeriodoEspacio"
$sheet.Range("G1").Value = "FiltroFecha"

# Data row
$sheet.Range("A2").Value = 271
$sheet.Range("B2").Value = "Octubre 30"
$sheet.Range("C2").Value = '"BANCO DE LA REPUBLICA NOMINA JUBILADOS"'
$sheet.Range("D2").Value = '"C:\Users\User\Downloads\PlanillasCarguePagaduria\"'
$sheet.Range("E2").Value = 2021
$sheet.Range("F2").Value = '"Octubre  30"'
$sheet.Range("G2").Value = "30/10/2021"

# Text-format the cells the same way the source workbook's table did
$sheet.Range("A2").NumberFormat = "@"
$sheet.Range("B1:G1").NumberFormat = "@"
$sheet.Range("B2:D2").NumberFormat = "@"
$sheet.Range("E2").NumberFormat = "@"
$sheet.Range("F2:G2").NumberFormat = "@"

# Column widths: B keeps the default table width, C/D best-fit to their
# (long) contents, same as the original workbook.
$sheet.Columns.Item(2).ColumnWidth = 11.42578125
$sheet.Columns.Item(3).AutoFit()
$sheet.Columns.Item(4).AutoFit()

# This is now the only/active sheet, selection parked on F2
[void]$sheet.Activate()
[void]$sheet.Range("F2").Select()
